# Commit "it's been a while": append new subjects S3 (nick), S4 (marijn), S5 (yelena)
# to the watchErpDataset sheet (rows 18-41), reproducing the original shared-string
# interning order: the label columns (sessionDirectory/subjectTag/subjectName) for the
# first row of each block are written first, then all fileName values for the S4 and S5
# blocks, then all fileName values for the S3 block -- matching how the source data was
# evidently appended to this lab log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First-row label cells, in original authoring order (creates new shared strings 28-36) ---
$ws.Range("D18").Value = "2013-06-12-nick"
$ws.Range("A18").Value = "S3"
$ws.Range("B18").Value = "nick"
$ws.Range("A26").Value = "S4"
$ws.Range("B26").Value = "marijn"
$ws.Range("D26").Value = "2013-06-14-marijn"
$ws.Range("A34").Value = "S5"
$ws.Range("B34").Value = "yelena"
$ws.Range("D34").Value = "2013-06-14-yelena"

# --- fileName values: S4 block (rows 26-33), then S5 block (rows 34-41) ---
$ws.Range("E26").Value = "2013-06-14-14-13-50-run1"
$ws.Range("E27").Value = "2013-06-14-14-21-21-run2"
$ws.Range("E28").Value = "2013-06-14-14-27-45-run3"
$ws.Range("E29").Value = "2013-06-14-14-33-18-run4"
$ws.Range("E30").Value = "2013-06-14-14-51-21-run5"
$ws.Range("E31").Value = "2013-06-14-14-58-54-run6"
$ws.Range("E32").Value = "2013-06-14-15-04-37-run7"
$ws.Range("E33").Value = "2013-06-14-15-11-23-run8"
$ws.Range("E34").Value = "2013-06-14-16-49-06-run1"
$ws.Range("E35").Value = "2013-06-14-16-55-50-run2"
$ws.Range("E36").Value = "2013-06-14-17-01-38-run3"
$ws.Range("E37").Value = "2013-06-14-17-07-22-run4"
$ws.Range("E38").Value = "2013-06-14-17-48-31-run5"
$ws.Range("E39").Value = "2013-06-14-17-55-46-run6"
$ws.Range("E40").Value = "2013-06-14-18-02-27-run7"
$ws.Range("E41").Value = "2013-06-14-18-10-22-run8"

# --- fileName values: S3 block (rows 18-25), written last ---
$ws.Range("E18").Value = "2013-06-12-14-55-14-run1"
$ws.Range("E19").Value = "2013-06-12-15-02-03-run2"
$ws.Range("E20").Value = "2013-06-12-15-07-46-run3"
$ws.Range("E21").Value = "2013-06-12-15-13-10-run4"
$ws.Range("E22").Value = "2013-06-12-15-24-23-run5"
$ws.Range("E23").Value = "2013-06-12-15-30-28-run6"
$ws.Range("E24").Value = "2013-06-12-15-39-52-run7"
$ws.Range("E25").Value = "2013-06-12-15-45-28-run8"

# --- Repeat the label text down each block (reuses the shared strings interned above) ---
$ws.Range("A19").Value = "S3"
$ws.Range("B19").Value = "nick"
$ws.Range("D19").Value = "2013-06-12-nick"
$ws.Range("A20").Value = "S3"
$ws.Range("B20").Value = "nick"
$ws.Range("D20").Value = "2013-06-12-nick"
$ws.Range("A21").Value = "S3"
$ws.Range("B21").Value = "nick"
$ws.Range("D21").Value = "2013-06-12-nick"
$ws.Range("A22").Value = "S3"
$ws.Range("B22").Value = "nick"
$ws.Range("D22").Value = "2013-06-12-nick"
$ws.Range("A23").Value = "S3"
$ws.Range("B23").Value = "nick"
$ws.Range("D23").Value = "2013-06-12-nick"
$ws.Range("A24").Value = "S3"
$ws.Range("B24").Value = "nick"
$ws.Range("D24").Value = "2013-06-12-nick"
$ws.Range("A25").Value = "S3"
$ws.Range("B25").Value = "nick"
$ws.Range("D25").Value = "2013-06-12-nick"
$ws.Range("A27").Value = "S4"
$ws.Range("B27").Value = "marijn"
$ws.Range("D27").Value = "2013-06-14-marijn"
$ws.Range("A28").Value = "S4"
$ws.Range("B28").Value = "marijn"
$ws.Range("D28").Value = "2013-06-14-marijn"
$ws.Range("A29").Value = "S4"
$ws.Range("B29").Value = "marijn"
$ws.Range("D29").Value = "2013-06-14-marijn"
$ws.Range("A30").Value = "S4"
$ws.Range("B30").Value = "marijn"
$ws.Range("D30").Value = "2013-06-14-marijn"
$ws.Range("A31").Value = "S4"
$ws.Range("B31").Value = "marijn"
$ws.Range("D31").Value = "2013-06-14-marijn"
$ws.Range("A32").Value = "S4"
$ws.Range("B32").Value = "marijn"
$ws.Range("D32").Value = "2013-06-14-marijn"
$ws.Range("A33").Value = "S4"
$ws.Range("B33").Value = "marijn"
$ws.Range("D33").Value = "2013-06-14-marijn"
$ws.Range("A35").Value = "S5"
$ws.Range("B35").Value = "yelena"
$ws.Range("D35").Value = "2013-06-14-yelena"
$ws.Range("A36").Value = "S5"
$ws.Range("B36").Value = "yelena"
$ws.Range("D36").Value = "2013-06-14-yelena"
$ws.Range("A37").Value = "S5"
$ws.Range("B37").Value = "yelena"
$ws.Range("D37").Value = "2013-06-14-yelena"
$ws.Range("A38").Value = "S5"
$ws.Range("B38").Value = "yelena"
$ws.Range("D38").Value = "2013-06-14-yelena"
$ws.Range("A39").Value = "S5"
$ws.Range("B39").Value = "yelena"
$ws.Range("D39").Value = "2013-06-14-yelena"
$ws.Range("A40").Value = "S5"
$ws.Range("B40").Value = "yelena"
$ws.Range("D40").Value = "2013-06-14-yelena"
$ws.Range("A41").Value = "S5"
$ws.Range("B41").Value = "yelena"
$ws.Range("D41").Value = "2013-06-14-yelena"

# --- Numeric columns: date serials (col C) and run numbers (col F) ---
$ws.Range("C18").Value = 41437
$ws.Range("C19").Value = 41437
$ws.Range("C20").Value = 41437
$ws.Range("C21").Value = 41437
$ws.Range("C22").Value = 41437
$ws.Range("C23").Value = 41437
$ws.Range("C24").Value = 41437
$ws.Range("C25").Value = 41437
$ws.Range("C26").Value = 41439
$ws.Range("C27").Value = 41439
$ws.Range("C28").Value = 41439
$ws.Range("C29").Value = 41439
$ws.Range("C30").Value = 41439
$ws.Range("C31").Value = 41439
$ws.Range("C32").Value = 41439
$ws.Range("C33").Value = 41439
$ws.Range("C34").Value = 41439
$ws.Range("C35").Value = 41439
$ws.Range("C36").Value = 41439
$ws.Range("C37").Value = 41439
$ws.Range("C38").Value = 41439
$ws.Range("C39").Value = 41439
$ws.Range("C40").Value = 41439
$ws.Range("C41").Value = 41439
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = 6
$ws.Range("F24").Value = 7
$ws.Range("F25").Value = 8
$ws.Range("F26").Value = 1
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 5
$ws.Range("F31").Value = 6
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 8
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 2
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 5
$ws.Range("F39").Value = 6
$ws.Range("F40").Value = 7
$ws.Range("F41").Value = 8

# --- Copy the formatting used by the existing S2 block (rows 10-17) onto all new rows ---
$ws.Range("A10:F10").Copy()
$ws.Range("A18:F41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Final selection state recorded in the sheet view ---
$ws.Range("F21").Select()
